$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): row 6 and row 11 "想去人数" (F column) counts updated
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 504
$wsExpo.Range("F11").Value = 696

# Sheet "全部类型" (all types): same two events duplicated at row 6 and row 15
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 504
$wsAll.Range("F15").Value = 696
